# guide41_movie.xlsx edit:
#  - Insert a new row above row 5 on sheet "p1", filled with a copy of the
#    B1 banner cell (value + format + row height).
#  - Make "p1" the active sheet/tab (was "p2").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("p1")

# Insert a blank row above row 5; existing rows 5.. shift down to 6..
$ws.Rows(5).Insert()

# Paste the banner cell's value (reuses the existing shared string) into
# the new B5.
$ws.Range("B1").Copy()
$ws.Range("B5").PasteSpecial(-4163)   # xlPasteValues

# Paste the banner cell's formatting onto B5 (new row) and B6 (the row
# that used to be row 5), matching the source workbook.
$ws.Range("B1").Copy()
$ws.Range("B5:B6").PasteSpecial(-4122)   # xlPasteFormats

# Match the banner row's auto-fit height on the newly inserted row.
$ws.Rows(5).RowHeight = 207.75

# Switch the active sheet/tab from "p2" to "p1".
$ws.Activate()
